$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "hy_w36": every row (B3:B84) previously held a unique label
#    hyline_w37 .. hyline_w118; they all get collapsed to "hyline_w36"
#    (same value already used in B2). Doing this through COM also makes
#    the now-unused shared strings disappear automatically, which is
#    what re-numbers every other <v> index referenced elsewhere in the
#    workbook (sheet5/sheet6/sheet7/sheet9 header cells, etc.).
# ---------------------------------------------------------------------
$wsHyW36 = $wb.Worksheets.Item("hy_w36")
for ($r = 3; $r -le 84; $r++) {
    $wsHyW36.Range("B$r").Value = "hyline_w36"
}

# ---------------------------------------------------------------------
# 2) Sheet "isa_brown": P2 used to reference a now-broken external
#    reference (=#REF!+O2) which propagated #REF! errors down the
#    whole shared-formula column. Point it back at O2 so the rest of
#    the shared formula chain (P3:P73, "=P(prev)+O(row)") recalculates
#    to real numbers again.
# ---------------------------------------------------------------------
$wsIsaBrown = $wb.Worksheets.Item("isa_brown")
$wsIsaBrown.Range("P2").Formula = "=O2"

# Column width tweaks on isa_brown (B, C, H)
$wsIsaBrown.Columns.Item(2).ColumnWidth = 14.627604166666666
$wsIsaBrown.Columns.Item(3).ColumnWidth = 10.549479166666666
$wsIsaBrown.Columns.Item(8).ColumnWidth = 20.858072916666668

# Update selection on isa_brown before it stops being the active sheet
$wsIsaBrown.Range("P15").Select()

# ---------------------------------------------------------------------
# 3) View state: hy_w36 becomes the active/selected tab (instead of
#    isa_brown), scrolled down and with a new selected cell.
# ---------------------------------------------------------------------
$wsHyW36.Activate()
$wsHyW36.Range("E75").Select()

$wb.Save()
